# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (F column) counts on the sheets that list individual
# records ("展览" and "全部类型"). Other sheets are left untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 8509
    $ws.Range("F4").Value = 6207
    $ws.Range("F5").Value = 542
    $ws.Range("F9").Value = 330

    if ($sheetName -eq "展览") {
        $ws.Range("F10").Value = 1179
    } elseif ($sheetName -eq "全部类型") {
        $ws.Range("F14").Value = 1179
    }
}
